# Update the "충돌체크" row of the progress table on slide 2:
#   - "실제 개발 완료 범위" cell: append ", 낙사" after "장애물과의 충돌"
#   - "진척도" cell: bump the progress value from 80% to 100%

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table

# Row 3 = "충돌체크"; Column 3 = "실제 개발 완료 범위"; Column 4 = "진척도"
$descCell = $tbl.Cell(3, 3)
$descRange = $descCell.Shape.TextFrame.TextRange
$descRange.Text = $descRange.Text + ", 낙사"

$pctCell = $tbl.Cell(3, 4)
$pctRange = $pctCell.Shape.TextFrame.TextRange
$pctRange.Text = "100%"
